# Applies a re-shuffling of several data rows in the "Artfynd" worksheet.
# Columns A, B, D, E, F, G, H, Q, R are relocated between rows according to
# the mapping below (derived from the target diff); all other columns are
# left untouched.
#
# target row -> source row (data is copied FROM source row's original values
# TO target row):
#   2  <- 12
#   3  <- 7
#   4  <- 8
#   6  <- 10
#   7  <- 2
#   8  <- 3
#   10 <- 6
#   11 <- 13
#   12 <- 4
#   13 <- 11
# (rows 5 and 9 are unaffected)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

$affectedRows = @(2, 3, 4, 6, 7, 8, 10, 11, 12, 13)

# Snapshot the original values of every affected row/column before writing
# anything, so source rows used later in the cycle aren't clobbered first.
$snapshot = @{}
foreach ($r in $affectedRows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

$mapping = @{
    2  = 12
    3  = 7
    4  = 8
    6  = 10
    7  = 2
    8  = 3
    10 = 6
    11 = 13
    12 = 4
    13 = 11
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $sourceVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value2 = $sourceVals[$c]
    }
}
